# Apply the "Updated cryptos list" refresh described by the commit diff.
# Only the cells whose text actually changed between before/after are touched;
# all other cells (column A ranks, and any untouched B/C/D/E cells) are left alone.
#
# The "Price" column (D) holds values that are display TEXT, not numbers
# (e.g. "69.400.19" uses dots as thousands separators, "7.20" has a
# significant trailing zero). Assigning such literals straight to .Value
# would let Excel auto-coerce the numeric-looking ones into real numbers,
# so for those specific cells we force the cell to Text format first -
# exactly like typing a leading apostrophe would do in the Excel UI -
# before writing the new value. Cells whose new text can never be parsed as
# a plain number (it still contains multiple "." separators) do not need this.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "69.400.19"
$ws.Range("E2").Value = "  -0.21%  "

# Row 3
$ws.Range("D3").Value = "3.663.33"
$ws.Range("E3").Value = "  -0.70%  "

# Row 4
$ws.Range("E4").Value = "  -0.09%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "623.47"
$ws.Range("E5").Value = "  -7.24%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.21"
$ws.Range("E6").Value = "  -1.06%  "

# Row 7
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.495"
$ws.Range("E8").Value = "  -0.51%  "

# Row 9
$ws.Range("E9").Value = "  -1.01%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.20"
$ws.Range("E10").Value = "  +0.94%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.439"
$ws.Range("E11").Value = "  -0.60%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000228"
$ws.Range("E12").Value = "  -2.31%  "

# Row 13
$ws.Range("D13").Value = "4.281.03"
$ws.Range("E13").Value = "  -0.86%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.18"
$ws.Range("E14").Value = "  -2.23%  "

# Row 15
$ws.Range("D15").Value = "3.686.43"
$ws.Range("E15").Value = "  +0.14%  "

# Row 16
$ws.Range("D16").Value = "69.447.87"
$ws.Range("E16").Value = "  -0.18%  "

# Row 17
$ws.Range("E17").Value = "  +1.29%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.49"
$ws.Range("E18").Value = "  -0.47%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "15.81"
$ws.Range("E19").Value = "  -2.34%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.32"
$ws.Range("E20").Value = "  +5.74%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "466.84"
$ws.Range("E21").Value = "  -1.38%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.648"
$ws.Range("E22").Value = "  -0.31%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "79.19"
$ws.Range("E23").Value = "  -1.43%  "

# Row 24
$ws.Range("D24").Value = "3.809.83"
$ws.Range("E24").Value = "  -0.77%  "

# Row 25
$ws.Range("E25").Value = "  +0.04%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.14"
$ws.Range("E26").Value = "  +1.32%  "

# Row 27
$ws.Range("E27").Value = "  -3.94%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.62"
$ws.Range("E28").Value = "  -6.06%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.61"
$ws.Range("E29").Value = "  -3.28%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.67"
$ws.Range("E30").Value = "  -3.18%  "

# Row 31
$ws.Range("E31").Value = "  -0.07%  "

# Row 32
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.96"
$ws.Range("E32").Value = "  -2.75%  "

# Row 33
$ws.Range("B33").Value = "Kaspa"
$ws.Range("C33").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.164"
$ws.Range("E33").Value = "  -2.51%  "

# Row 34
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.49"
$ws.Range("E34").Value = "  -1.43%  "

# Row 35
$ws.Range("D35").Value = "3.666.86"
$ws.Range("E35").Value = "  -0.54%  "

# Row 36
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.36"
$ws.Range("E36").Value = "  -2.42%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.27"
$ws.Range("E37").Value = "  -2.91%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "179.63"
$ws.Range("E39").Value = "  +2.64%  "

# Row 40
$ws.Range("E40").Value = "  -0.13%  "

# Row 41
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.80"
$ws.Range("E41").Value = "  -5.34%  "

# Row 42
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.22"
$ws.Range("E42").Value = "  -1.38%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0895"
$ws.Range("E43").Value = "  -1.46%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.921"
$ws.Range("E44").Value = "  -1.80%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "29.49"
$ws.Range("E45").Value = "  +6.66%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "46.70"
$ws.Range("E46").Value = "  -0.51%  "

# Row 47
$ws.Range("E47").Value = "  -1.91%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.83"
$ws.Range("E48").Value = "  -0.52%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000262"
$ws.Range("E49").Value = "  -5.00%  "

# Row 50
$ws.Range("E50").Value = "  -5.02%  "

# Row 51
$ws.Range("B51").Value = "ONDO"
$ws.Range("C51").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.20"
$ws.Range("E51").Value = "  -7.03%  "
